# Daily attendance processing - 2026-01-28 11:39:45
# Reorders the "Recorded By" (column G) list of names/emails on the
# "Session Analysis Results" sheet so that email addresses come first,
# followed by the remaining (non-email) entries, preserving the
# relative order within each group.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $value = $cell.Value2

    if ($value -ne $null -and $value -ne "") {
        $parts = $value -split ", "

        if ($parts.Count -gt 1) {
            $emails = @()
            $others = @()

            foreach ($part in $parts) {
                if ($part -like "*@*") {
                    $emails += $part
                } else {
                    $others += $part
                }
            }

            $reordered = $emails + $others
            $newValue = $reordered -join ", "

            if ($newValue -ne $value) {
                $cell.Value = $newValue
            }
        }
    }
}
